# generate data set api added
# Appends newly-generated simulation rows to the three sheets:
#   ReachIntersection   -> rows 9-16  (A:D)
#   SpeedAtIntersection -> rows 15-25 (A:C)
#   SetFinalSpeed       -> row 8      (A:D)
#
# All source values in this workbook are stored as literal text (not
# numbers), so every cell is force-formatted as Text ("@") before its
# value is assigned - this preserves exact string representations such
# as trailing ".0" and long floating point literals.

$wb = $excel.ActiveWorkbook

# NOTE: this COM engine only binds *positional* parameters reliably, so
# the helper below is called without the `-Name value` syntax.
function Set-TextRow {
    param($Sheet, $Row, $Values)
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $cell = $Sheet.Cells.Item($Row, $i + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $Values[$i]
    }
}

# ---------------------------------------------------------------------
# Sheet 1: ReachIntersection  (columns: light, distance, speed, direction)
# ---------------------------------------------------------------------
$wsReach = $wb.Worksheets.Item("ReachIntersection")

Set-TextRow $wsReach 9  @("1", "1165.492770458883", "-666.0", "740")
Set-TextRow $wsReach 10 @("0", "1165.492770458883", "-1687.7", "-5271")
Set-TextRow $wsReach 11 @("0", "1609.0", "4060.344", "0")
Set-TextRow $wsReach 12 @("2", "956.795", "133.57108176521356", "1730")
Set-TextRow $wsReach 13 @("2", "956.795", "2174.030431", "4388")
Set-TextRow $wsReach 14 @("1", "956.795", "956.795", "4388")
Set-TextRow $wsReach 15 @("0", "-1239.0", "-1123.5137", "-680")
Set-TextRow $wsReach 16 @("0", "133.5711", "847.418", "0")

# ---------------------------------------------------------------------
# Sheet 2: SpeedAtIntersection  (columns: accelFlag, speed, distance)
# ---------------------------------------------------------------------
$wsSpeed = $wb.Worksheets.Item("SpeedAtIntersection")

Set-TextRow $wsSpeed 15 @("-1900", "4544.8", "4544.8")
Set-TextRow $wsSpeed 16 @("-1900", "-3119.4", "-4406.942272")
Set-TextRow $wsSpeed 17 @("1", "1227.23", "1227.23")
Set-TextRow $wsSpeed 18 @("1", "1227.23", "1227.23")
Set-TextRow $wsSpeed 19 @("-853", "526.0", "240.9")
Set-TextRow $wsSpeed 20 @("-2261", "526.0", "-2255.709063")
Set-TextRow $wsSpeed 21 @("2419", "228.201385", "1501.4845")
Set-TextRow $wsSpeed 22 @("-70", "956.795", "-663.4069")
Set-TextRow $wsSpeed 23 @("-70", "-1518.0", "-1518.0")
Set-TextRow $wsSpeed 24 @("-70", "1076.1574168684874", "1076.1574168684874")
Set-TextRow $wsSpeed 25 @("-70", "1432.45155", "1432.45155")

# ---------------------------------------------------------------------
# Sheet 3: SetFinalSpeed  (columns: oldSpeed, speed, lowLimit, highLimit)
# ---------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("SetFinalSpeed")

Set-TextRow $wsFinal 8 @("1947.180595", "1947.180595", "-571", "-259")

Write-Output "generated data set rows appended"
